$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the resolution comment for row 8 (D8) to "Incorporated.", matching
# the same value already used in D3:D7.
$ws.Range("D8").Value = "Incorporated."

# Move the active selection to D9 (next cell down).
$ws.Range("D9").Select()
